# Update the Build/TaskList "Sheet1" task-progress column (C) to reflect
# current build & task status, per commit "The build and task progress".
#
#   Row  8 (task 7)  : In progress -> Completed
#   Row 10 (task 9)  : Not started -> Completed
#   Row 11 (task 10) : Not started -> Completed
#   Row 12 (task 11) : Not started -> Completed
#   Row 13 (task 12) : Not started -> In progress
#   Row 14 (task 13) : Not started -> Completed
#   Row 16 (task 15) : Not started -> Completed
#
# (Rows 2-7, 9 and 15 keep their existing status.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C8").Value = "Completed"
$ws.Range("C10").Value = "Completed"
$ws.Range("C11").Value = "Completed"
$ws.Range("C12").Value = "Completed"
$ws.Range("C13").Value = "In progress"
$ws.Range("C14").Value = "Completed"
$ws.Range("C16").Value = "Completed"

# Leave the view scrolled to the top and the selection where the user's
# last edit landed (C18, just past the last data row).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C18").Select()
